# Natmi following Dr Hou advice
# Recompute the Thbs2-Itga6 LR-pairs table: the sending/target cluster cross
# join now also includes the "ECs" cluster (in addition to FAPs and sCs),
# giving a full 3x3 grid (rows 2-10), and all numeric metrics are refreshed
# to match the recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Thbs2"
$ws.Cells.Item(2, 3).Value = "Itga6"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.8911683333333333
$ws.Cells.Item(2, 8).Value = 2.673505
$ws.Cells.Item(2, 9).Value = 0.02693425114262819
$ws.Cells.Item(2, 10).Value = 0.02693425114262819
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 147.4213356666667
$ws.Cells.Item(2, 14).Value = 442.264007
$ws.Cells.Item(2, 15).Value = 0.9507885170992249
$ws.Cells.Item(2, 16).Value = 0.950788517099225
$ws.Cells.Item(2, 17).Value = 131.3772260038372
$ws.Cells.Item(2, 18).Value = 1182.395034034535
$ws.Cells.Item(2, 19).Value = 0.02560877670307755
$ws.Cells.Item(2, 20).Value = 0.02560877670307756

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Thbs2"
$ws.Cells.Item(3, 3).Value = "Itga6"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.8911683333333333
$ws.Cells.Item(3, 8).Value = 2.673505
$ws.Cells.Item(3, 9).Value = 0.02693425114262819
$ws.Cells.Item(3, 10).Value = 0.02693425114262819
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.340788333333334
$ws.Cells.Item(3, 14).Value = 7.022365000000001
$ws.Cells.Item(3, 15).Value = 0.01509682881537204
$ws.Cells.Item(3, 16).Value = 0.01509682881537204
$ws.Cells.Item(3, 17).Value = 2.086036437702778
$ws.Cells.Item(3, 18).Value = 18.774327939325
$ws.Cells.Item(3, 19).Value = 0.0004066217787704964
$ws.Cells.Item(3, 20).Value = 0.0004066217787704964

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Thbs2"
$ws.Cells.Item(4, 3).Value = "Itga6"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.8911683333333333
$ws.Cells.Item(4, 8).Value = 2.673505
$ws.Cells.Item(4, 9).Value = 0.02693425114262819
$ws.Cells.Item(4, 10).Value = 0.02693425114262819
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 5.289533666666667
$ws.Cells.Item(4, 14).Value = 15.868601
$ws.Cells.Item(4, 15).Value = 0.03411465408540306
$ws.Cells.Item(4, 16).Value = 0.03411465408540307
$ws.Cells.Item(4, 17).Value = 4.713864901833889
$ws.Cells.Item(4, 18).Value = 42.424784116505
$ws.Cells.Item(4, 19).Value = 0.0009188526607801327
$ws.Cells.Item(4, 20).Value = 0.0009188526607801329

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Thbs2"
$ws.Cells.Item(5, 3).Value = "Itga6"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 24.359699
$ws.Cells.Item(5, 8).Value = 73.07909699999999
$ws.Cells.Item(5, 9).Value = 0.7362360466408275
$ws.Cells.Item(5, 10).Value = 0.7362360466408276
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 147.4213356666667
$ws.Cells.Item(5, 14).Value = 442.264007
$ws.Cells.Item(5, 15).Value = 0.9507885170992249
$ws.Cells.Item(5, 16).Value = 0.950788517099225
$ws.Cells.Item(5, 17).Value = 3591.139363017964
$ws.Cells.Item(5, 18).Value = 32320.25426716168
$ws.Cells.Item(5, 19).Value = 0.7000047790206282
$ws.Cells.Item(5, 20).Value = 0.7000047790206283

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Thbs2"
$ws.Cells.Item(6, 3).Value = "Itga6"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 24.359699
$ws.Cells.Item(6, 8).Value = 73.07909699999999
$ws.Cells.Item(6, 9).Value = 0.7362360466408275
$ws.Cells.Item(6, 10).Value = 0.7362360466408276
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.340788333333334
$ws.Cells.Item(6, 14).Value = 7.022365000000001
$ws.Cells.Item(6, 15).Value = 0.01509682881537204
$ws.Cells.Item(6, 16).Value = 0.01509682881537204
$ws.Cells.Item(6, 17).Value = 57.02089922271166
$ws.Cells.Item(6, 18).Value = 513.188093004405
$ws.Cells.Item(6, 19).Value = 0.01111482956384283
$ws.Cells.Item(6, 20).Value = 0.01111482956384284

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Thbs2"
$ws.Cells.Item(7, 3).Value = "Itga6"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 24.359699
$ws.Cells.Item(7, 8).Value = 73.07909699999999
$ws.Cells.Item(7, 9).Value = 0.7362360466408275
$ws.Cells.Item(7, 10).Value = 0.7362360466408276
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 5.289533666666667
$ws.Cells.Item(7, 14).Value = 15.868601
$ws.Cells.Item(7, 15).Value = 0.03411465408540306
$ws.Cells.Item(7, 16).Value = 0.03411465408540307
$ws.Cells.Item(7, 17).Value = 128.8514479703663
$ws.Cells.Item(7, 18).Value = 1159.663031733297
$ws.Cells.Item(7, 19).Value = 0.0251164380563565
$ws.Cells.Item(7, 20).Value = 0.02511643805635651

$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Thbs2"
$ws.Cells.Item(8, 3).Value = "Itga6"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 7.835938333333334
$ws.Cells.Item(8, 8).Value = 23.507815
$ws.Cells.Item(8, 9).Value = 0.2368297022165442
$ws.Cells.Item(8, 10).Value = 0.2368297022165442
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 147.4213356666667
$ws.Cells.Item(8, 14).Value = 442.264007
$ws.Cells.Item(8, 15).Value = 0.9507885170992249
$ws.Cells.Item(8, 16).Value = 0.950788517099225
$ws.Cells.Item(8, 17).Value = 1155.184495301634
$ws.Cells.Item(8, 18).Value = 10396.6604577147
$ws.Cells.Item(8, 19).Value = 0.2251749613755191
$ws.Cells.Item(8, 20).Value = 0.2251749613755191

$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Thbs2"
$ws.Cells.Item(9, 3).Value = "Itga6"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 7.835938333333334
$ws.Cells.Item(9, 8).Value = 23.507815
$ws.Cells.Item(9, 9).Value = 0.2368297022165442
$ws.Cells.Item(9, 10).Value = 0.2368297022165442
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.340788333333334
$ws.Cells.Item(9, 14).Value = 7.022365000000001
$ws.Cells.Item(9, 15).Value = 0.01509682881537204
$ws.Cells.Item(9, 16).Value = 0.01509682881537204
$ws.Cells.Item(9, 17).Value = 18.34227303138611
$ws.Cells.Item(9, 18).Value = 165.080457282475
$ws.Cells.Item(9, 19).Value = 0.003575377472758703
$ws.Cells.Item(9, 20).Value = 0.003575377472758703

$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Thbs2"
$ws.Cells.Item(10, 3).Value = "Itga6"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 7.835938333333334
$ws.Cells.Item(10, 8).Value = 23.507815
$ws.Cells.Item(10, 9).Value = 0.2368297022165442
$ws.Cells.Item(10, 10).Value = 0.2368297022165442
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 5.289533666666667
$ws.Cells.Item(10, 14).Value = 15.868601
$ws.Cells.Item(10, 15).Value = 0.03411465408540306
$ws.Cells.Item(10, 16).Value = 0.03411465408540307
$ws.Cells.Item(10, 17).Value = 41.44845962409056
$ws.Cells.Item(10, 18).Value = 373.036136616815
$ws.Cells.Item(10, 19).Value = 0.008079363368266422
$ws.Cells.Item(10, 20).Value = 0.008079363368266423
